# Daily attendance processing - 2025-10-08 08:50:00
# Applies the "Recorded By" ordering fixes, updated session statistics,
# and marks the B2D / B2E / B2F Session-14 rows as Recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper scratch cell used to write literal percentage TEXT (e.g. "53.8%")
# without Excel's automatic "looks like a percentage" numeric conversion.
# Formatting the scratch cell as Text ("@") first, writing the string,
# then copying *values only* back onto the real target preserves the
# target cell's existing style while keeping the content as plain text.
# ---------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

# ---------------------------------------------------------------------
# "Recorded By" (column G) name-order fixes
# ---------------------------------------------------------------------
$ws.Range("G2").Value   = "backup@backdoor.com, System, system"
$ws.Range("G4").Value   = "backup@backdoor.com, System"
$ws.Range("G5").Value   = "backup@backdoor.com, System"
$ws.Range("G11").Value  = "System, dnasr281@gmail.com"
$ws.Range("G29").Value  = "backup@backdoor.com, System, system"
$ws.Range("G32").Value  = "backup@backdoor.com, System"
$ws.Range("G38").Value  = "System, dnasr281@gmail.com"
$ws.Range("G56").Value  = "backup@backdoor.com, System, system"
$ws.Range("G58").Value  = "backup@backdoor.com, System"
$ws.Range("G59").Value  = "backup@backdoor.com, System"
$ws.Range("G65").Value  = "System, dnasr281@gmail.com"
$ws.Range("G84").Value  = "backup@backdoor.com, System"
$ws.Range("G85").Value  = "backup@backdoor.com, System"
$ws.Range("G90").Value  = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G110").Value = "backup@backdoor.com, System"
$ws.Range("G111").Value = "backup@backdoor.com, System"
$ws.Range("G116").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G136").Value = "backup@backdoor.com, System"
$ws.Range("G137").Value = "backup@backdoor.com, System"
$ws.Range("G142").Value = "admin@admin.com, dnasr281@gmail.com"

# ---------------------------------------------------------------------
# Class statistics block (K6:L10)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 86
$ws.Range("L8").Value = 72

$scratch.Value = "54.1%"
$scratch.Copy()
$ws.Range("L9").PasteSpecial(-4163)

$scratch.Value = "70.3%"
$scratch.Copy()
$ws.Range("L10").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# Per-group summary table (rows 18-20 => groups B2D, B2E, B2F)
# ---------------------------------------------------------------------
$ws.Range("O18").Value = 14
$ws.Range("Q18").Value = 12
$scratch.Value = "53.8%"
$scratch.Copy()
$ws.Range("R18").PasteSpecial(-4163)

$ws.Range("O19").Value = 14
$ws.Range("Q19").Value = 12
$scratch.Value = "53.8%"
$scratch.Copy()
$ws.Range("R19").PasteSpecial(-4163)
$scratch.Value = "75.8%"
$scratch.Copy()
$ws.Range("S19").PasteSpecial(-4163)

$ws.Range("O20").Value = 14
$ws.Range("Q20").Value = 12
$scratch.Value = "53.8%"
$scratch.Copy()
$ws.Range("R20").PasteSpecial(-4163)
$scratch.Value = "74.7%"
$scratch.Copy()
$ws.Range("S20").PasteSpecial(-4163)

$scratch.Clear()

# ---------------------------------------------------------------------
# Session 14 (08/10/2025) rows for B2D / B2E / B2F flip from
# "Pending" (yellow, style 5) to "Recorded" (green, style 2).
#
# Strategy: copy a same-shaped already-"Recorded" row (same style 2)
# over the target row to pick up both content + the green format in one
# shot, then restore the Session/Date columns (D/E) - which must stay
# exactly as they were - using a values-only paste so the newly-applied
# green format is not disturbed, and finally set the real
# Recorded-By / Students / Status values.
# ---------------------------------------------------------------------

# --- Row 96 (B2D, Session 14) -----------------------------------------
$ws.Range("D96:E96").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4163)

$ws.Range("A90:I90").Copy()
$ws.Range("A96:I96").PasteSpecial(-4122)

$ws.Range("AA1:AB1").Copy()
$ws.Range("D96:E96").PasteSpecial(-4163)
$ws.Range("AA1:AB1").Clear()

$ws.Range("G96").Value = "dnasr281@gmail.com"
$ws.Range("H96").Value = "42/56"
$ws.Range("I96").Value = "Recorded"

# --- Row 122 (B2E, Session 14) -----------------------------------------
$ws.Range("D122:E122").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4163)

$ws.Range("A116:I116").Copy()
$ws.Range("A122:I122").PasteSpecial(-4122)

$ws.Range("AA1:AB1").Copy()
$ws.Range("D122:E122").PasteSpecial(-4163)
$ws.Range("AA1:AB1").Clear()

$ws.Range("G122").Value = "dnasr281@gmail.com"
$ws.Range("H122").Value = "40/55"
$ws.Range("I122").Value = "Recorded"

# --- Row 148 (B2F, Session 14) -----------------------------------------
$ws.Range("D148:E148").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4163)

$ws.Range("A142:I142").Copy()
$ws.Range("A148:I148").PasteSpecial(-4122)

$ws.Range("AA1:AB1").Copy()
$ws.Range("D148:E148").PasteSpecial(-4163)
$ws.Range("AA1:AB1").Clear()

$ws.Range("G148").Value = "dnasr281@gmail.com"
$ws.Range("H148").Value = "48/57"
$ws.Range("I148").Value = "Recorded"
